$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "sub-header" cells in row 3 under the date-range columns (C, D, E)
# now carry a "date" label, centered horizontally, while keeping their
# existing left/right border and top-vertical alignment.
$dateLabels = $ws.Range("C3:E3")
$dateLabels.Value = "date"
$dateLabels.HorizontalAlignment = -4108   # xlCenter
$dateLabels.VerticalAlignment = -4160     # xlTop

# Column G (the last date column) now also uses the text number format
# ("@"), matching column F's existing text format so both render the
# same way.
$ws.Range("G3").NumberFormat = "@"
